# "Flow of the project established"
# Update the Continuous Variables (after removing outliers) summary table
# with the recalculated descriptive statistics, tidy up column G's width
# so it matches column F, and leave the selection where it was left after
# the last edit (I12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- HRLYEARN row (row 2): Mean, Median, Min, Max, Standard Deviation ---
$ws.Range("D2").Value = 3.53
$ws.Range("E2").Value = 3.51
$ws.Range("F2").Value = 1.93
$ws.Range("G2").Value = 5.33
$ws.Range("H2").Value = 0.44

# --- UHRSMAIN row (row 3) ---
$ws.Range("D3").Value = 38.77
$ws.Range("F3").Value = 35
$ws.Range("H3").Value = 1.87

# --- TENURE row (row 4) ---
$ws.Range("D4").Value = 100.89
$ws.Range("H4").Value = 83.02

# Column G was slightly wider than column F; align it with F's (best-fit)
# width now that the values in both columns are comparable in size.
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Move/leave the active selection on I12, matching where the author's
# cursor ended up after finishing these edits.
$ws.Range("I12").Select()
